# Apply the August bioassay (#4) re-dating/relabeling edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6: original T0 timepoint dated 45455 -> 45525 (August)
# Rows 7-46: original timepoint dated 45458 -> 45528 (August)
# Column A (Bioassay) changes from 1 to 4 for all data rows.
# Column B (Month) changes from "May" to "August" for all data rows.
# Columns D, E, F are left unchanged.

for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 1).Value2 = 4
    $ws.Cells.Item($r, 2).Value2 = "August"
    if ($r -le 6) {
        $ws.Cells.Item($r, 3).Value2 = 45525
    } else {
        $ws.Cells.Item($r, 3).Value2 = 45528
    }
}

# Update the sheet selection to match the post-edit state (column C selected).
[void]$ws.Columns("C").Select()
